function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
Set-TextValue $ws.Cells.Item(256, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(256, 2) "14:22:34"
Set-TextValue $ws.Cells.Item(256, 3) "14:00"
Set-TextValue $ws.Cells.Item(256, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(256, 5) "No Motion"
Set-TextValue $ws.Cells.Item(256, 6) "Inactive"
Set-TextValue $ws.Cells.Item(257, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(257, 2) "14:22:36"
Set-TextValue $ws.Cells.Item(257, 3) "14:00"
Set-TextValue $ws.Cells.Item(257, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(257, 5) "No Motion"
Set-TextValue $ws.Cells.Item(257, 6) "Inactive"
Set-TextValue $ws.Cells.Item(258, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(258, 2) "14:22:38"
Set-TextValue $ws.Cells.Item(258, 3) "14:00"
Set-TextValue $ws.Cells.Item(258, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(258, 5) "No Motion"
Set-TextValue $ws.Cells.Item(258, 6) "Inactive"
Set-TextValue $ws.Cells.Item(259, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(259, 2) "14:22:43"
Set-TextValue $ws.Cells.Item(259, 3) "14:00"
Set-TextValue $ws.Cells.Item(259, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(259, 5) "No Motion"
Set-TextValue $ws.Cells.Item(259, 6) "Inactive"
Set-TextValue $ws.Cells.Item(260, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(260, 2) "14:22:48"
Set-TextValue $ws.Cells.Item(260, 3) "14:00"
Set-TextValue $ws.Cells.Item(260, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(260, 5) "No Motion"
Set-TextValue $ws.Cells.Item(260, 6) "Inactive"
Set-TextValue $ws.Cells.Item(261, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(261, 2) "14:22:51"
Set-TextValue $ws.Cells.Item(261, 3) "14:00"
Set-TextValue $ws.Cells.Item(261, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(261, 5) "Motion Detected"
Set-TextValue $ws.Cells.Item(261, 6) "Active"
Set-TextValue $ws.Cells.Item(262, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(262, 2) "14:22:58"
Set-TextValue $ws.Cells.Item(262, 3) "14:00"
Set-TextValue $ws.Cells.Item(262, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(262, 5) "No Motion"
Set-TextValue $ws.Cells.Item(262, 6) "Inactive"
Set-TextValue $ws.Cells.Item(263, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(263, 2) "14:23:03"
Set-TextValue $ws.Cells.Item(263, 3) "14:00"
Set-TextValue $ws.Cells.Item(263, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(263, 5) "No Motion"
Set-TextValue $ws.Cells.Item(263, 6) "Inactive"
Set-TextValue $ws.Cells.Item(264, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(264, 2) "14:23:04"
Set-TextValue $ws.Cells.Item(264, 3) "14:00"
Set-TextValue $ws.Cells.Item(264, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(264, 5) "Motion Detected"
Set-TextValue $ws.Cells.Item(264, 6) "Active"
Set-TextValue $ws.Cells.Item(265, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(265, 2) "14:23:11"
Set-TextValue $ws.Cells.Item(265, 3) "14:00"
Set-TextValue $ws.Cells.Item(265, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(265, 5) "No Motion"
Set-TextValue $ws.Cells.Item(265, 6) "Inactive"
Set-TextValue $ws.Cells.Item(266, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(266, 2) "14:23:16"
Set-TextValue $ws.Cells.Item(266, 3) "14:00"
Set-TextValue $ws.Cells.Item(266, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(266, 5) "No Motion"
Set-TextValue $ws.Cells.Item(266, 6) "Inactive"
Set-TextValue $ws.Cells.Item(267, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(267, 2) "14:23:21"
Set-TextValue $ws.Cells.Item(267, 3) "14:00"
Set-TextValue $ws.Cells.Item(267, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(267, 5) "No Motion"
Set-TextValue $ws.Cells.Item(267, 6) "Inactive"
Set-TextValue $ws.Cells.Item(268, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(268, 2) "14:23:26"
Set-TextValue $ws.Cells.Item(268, 3) "14:00"
Set-TextValue $ws.Cells.Item(268, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(268, 5) "No Motion"
Set-TextValue $ws.Cells.Item(268, 6) "Inactive"
Set-TextValue $ws.Cells.Item(269, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(269, 2) "14:23:31"
Set-TextValue $ws.Cells.Item(269, 3) "14:00"
Set-TextValue $ws.Cells.Item(269, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(269, 5) "No Motion"
Set-TextValue $ws.Cells.Item(269, 6) "Inactive"

$ws = $wb.Worksheets.Item("Humidity")
Set-TextValue $ws.Cells.Item(213, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(213, 2) "14:22:34"
Set-TextValue $ws.Cells.Item(213, 3) "14:00"
Set-TextValue $ws.Cells.Item(213, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(213, 5) "77.5%"
Set-TextValue $ws.Cells.Item(213, 6) "Active"
Set-TextValue $ws.Cells.Item(214, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(214, 2) "14:22:36"
Set-TextValue $ws.Cells.Item(214, 3) "14:00"
Set-TextValue $ws.Cells.Item(214, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(214, 5) "78.2%"
Set-TextValue $ws.Cells.Item(214, 6) "Active"
Set-TextValue $ws.Cells.Item(215, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(215, 2) "14:22:39"
Set-TextValue $ws.Cells.Item(215, 3) "14:00"
Set-TextValue $ws.Cells.Item(215, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(215, 5) "77.2%"
Set-TextValue $ws.Cells.Item(215, 6) "Active"
Set-TextValue $ws.Cells.Item(216, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(216, 2) "14:22:44"
Set-TextValue $ws.Cells.Item(216, 3) "14:00"
Set-TextValue $ws.Cells.Item(216, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(216, 5) "78.1%"
Set-TextValue $ws.Cells.Item(216, 6) "Active"
Set-TextValue $ws.Cells.Item(217, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(217, 2) "14:22:49"
Set-TextValue $ws.Cells.Item(217, 3) "14:00"
Set-TextValue $ws.Cells.Item(217, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(217, 5) "77.4%"
Set-TextValue $ws.Cells.Item(217, 6) "Active"
Set-TextValue $ws.Cells.Item(218, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(218, 2) "14:22:59"
Set-TextValue $ws.Cells.Item(218, 3) "14:00"
Set-TextValue $ws.Cells.Item(218, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(218, 5) "77.3%"
Set-TextValue $ws.Cells.Item(218, 6) "Active"
Set-TextValue $ws.Cells.Item(219, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(219, 2) "14:23:04"
Set-TextValue $ws.Cells.Item(219, 3) "14:00"
Set-TextValue $ws.Cells.Item(219, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(219, 5) "78.5%"
Set-TextValue $ws.Cells.Item(219, 6) "Active"
Set-TextValue $ws.Cells.Item(220, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(220, 2) "14:23:15"
Set-TextValue $ws.Cells.Item(220, 3) "14:00"
Set-TextValue $ws.Cells.Item(220, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(220, 5) "78.5%"
Set-TextValue $ws.Cells.Item(220, 6) "Active"
Set-TextValue $ws.Cells.Item(221, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(221, 2) "14:23:19"
Set-TextValue $ws.Cells.Item(221, 3) "14:00"
Set-TextValue $ws.Cells.Item(221, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(221, 5) "77.6%"
Set-TextValue $ws.Cells.Item(221, 6) "Active"
Set-TextValue $ws.Cells.Item(222, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(222, 2) "14:23:25"
Set-TextValue $ws.Cells.Item(222, 3) "14:00"
Set-TextValue $ws.Cells.Item(222, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(222, 5) "78.6%"
Set-TextValue $ws.Cells.Item(222, 6) "Active"
Set-TextValue $ws.Cells.Item(223, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(223, 2) "14:23:29"
Set-TextValue $ws.Cells.Item(223, 3) "14:00"
Set-TextValue $ws.Cells.Item(223, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(223, 5) "77.8%"
Set-TextValue $ws.Cells.Item(223, 6) "Active"

$ws = $wb.Worksheets.Item("Temperature")
Set-TextValue $ws.Cells.Item(213, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(213, 2) "14:22:35"
Set-TextValue $ws.Cells.Item(213, 3) "14:00"
Set-TextValue $ws.Cells.Item(213, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(213, 5) "24.7C"
Set-TextValue $ws.Cells.Item(213, 6) "Active"
Set-TextValue $ws.Cells.Item(214, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(214, 2) "14:22:37"
Set-TextValue $ws.Cells.Item(214, 3) "14:00"
Set-TextValue $ws.Cells.Item(214, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(214, 5) "24.7C"
Set-TextValue $ws.Cells.Item(214, 6) "Active"
Set-TextValue $ws.Cells.Item(215, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(215, 2) "14:22:40"
Set-TextValue $ws.Cells.Item(215, 3) "14:00"
Set-TextValue $ws.Cells.Item(215, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(215, 5) "24.7C"
Set-TextValue $ws.Cells.Item(215, 6) "Active"
Set-TextValue $ws.Cells.Item(216, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(216, 2) "14:22:45"
Set-TextValue $ws.Cells.Item(216, 3) "14:00"
Set-TextValue $ws.Cells.Item(216, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(216, 5) "24.7C"
Set-TextValue $ws.Cells.Item(216, 6) "Active"
Set-TextValue $ws.Cells.Item(217, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(217, 2) "14:22:50"
Set-TextValue $ws.Cells.Item(217, 3) "14:00"
Set-TextValue $ws.Cells.Item(217, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(217, 5) "24.7C"
Set-TextValue $ws.Cells.Item(217, 6) "Active"
Set-TextValue $ws.Cells.Item(218, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(218, 2) "14:23:00"
Set-TextValue $ws.Cells.Item(218, 3) "14:00"
Set-TextValue $ws.Cells.Item(218, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(218, 5) "24.6C"
Set-TextValue $ws.Cells.Item(218, 6) "Active"
Set-TextValue $ws.Cells.Item(219, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(219, 2) "14:23:05"
Set-TextValue $ws.Cells.Item(219, 3) "14:00"
Set-TextValue $ws.Cells.Item(219, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(219, 5) "24.7C"
Set-TextValue $ws.Cells.Item(219, 6) "Active"
Set-TextValue $ws.Cells.Item(220, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(220, 2) "14:23:15"
Set-TextValue $ws.Cells.Item(220, 3) "14:00"
Set-TextValue $ws.Cells.Item(220, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(220, 5) "24.6C"
Set-TextValue $ws.Cells.Item(220, 6) "Active"
Set-TextValue $ws.Cells.Item(221, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(221, 2) "14:23:20"
Set-TextValue $ws.Cells.Item(221, 3) "14:00"
Set-TextValue $ws.Cells.Item(221, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(221, 5) "24.6C"
Set-TextValue $ws.Cells.Item(221, 6) "Active"
Set-TextValue $ws.Cells.Item(222, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(222, 2) "14:23:25"
Set-TextValue $ws.Cells.Item(222, 3) "14:00"
Set-TextValue $ws.Cells.Item(222, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(222, 5) "24.6C"
Set-TextValue $ws.Cells.Item(222, 6) "Active"
Set-TextValue $ws.Cells.Item(223, 1) "2026-02-04"
Set-TextValue $ws.Cells.Item(223, 2) "14:23:30"
Set-TextValue $ws.Cells.Item(223, 3) "14:00"
Set-TextValue $ws.Cells.Item(223, 4) "Bathroom"
Set-TextValue $ws.Cells.Item(223, 5) "24.5C"
Set-TextValue $ws.Cells.Item(223, 6) "Active"
